# Append the next day's GSC export row to the "Chart" sheet:
#   Date = 2025-12-04, Non-HTTPS URLs = 0, HTTPS URLs = 26
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 60
$dateCell = $ws.Range("A" + $lastRow)

# Write the date as literal text (prefixing with an apostrophe stops Excel's
# auto-detection turning it into a date serial number), then re-paste the
# format from the previous row so the new cell ends up with the same
# (default) style as the rest of the column instead of a quote-prefixed one.
$dateCell.Value = "'2025-12-04"
$ws.Range("A59").Copy()
$dateCell.PasteSpecial(-4122)

$ws.Range("B" + $lastRow).Value = 0
$ws.Range("C" + $lastRow).Value = 26
